# M1 -> COMPLETED (add android libraries)
#
# Recolors several "label:" runs from the theme accent6/lumMod75% scheme
# color to an explicit red (FF0000), and folds the leading/trailing
# colon-space punctuation that used to live in a separate run into the
# (now red) label run, on slides 6, 8 and 9.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 6 ("Android Architecture") - Content Placeholder 2 (shape id 3)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

# "APPS" / ": where is located preinstalled apps and third part apps"
$tr6.Runs(2).Text = "APPS: "
$tr6.Runs(3).Text = "where is located preinstalled apps and third part apps"

# "ANDROID FRAMEWORK" / ": When you develop ..."
$tr6.Runs(4).Text = "ANDROID FRAMEWORK:"
$tr6.Runs(4).Font.Color.RGB = 255
$tr6.Runs(5).Text = " When you develop android app you use Android Framework, which offer an interface for native languages libraries (Java, C++)."

# "ANDROID RUNTIME" / ": is the android virtual machine ..."
$tr6.Runs(6).Text = "ANDROID RUNTIME: "
$tr6.Runs(6).Font.Color.RGB = 255
$tr6.Runs(7).Text = "is the android " + [char]0x201C + "virtual machine" + [char]0x201D + ". Core Libraries is a set of classes which are used by SDK and included functionality near the android VM"

# "HAL" (text unchanged, just recolor)
$tr6.Runs(8).Font.Color.RGB = 255

# "LINUX KERNEL" / ": " -> merge the colon-space into the label run and
# recolor; the now-empty run 11 is removed last so earlier indices stay
# valid.
$tr6.Runs(10).Text = "LINUX KERNEL: "
$tr6.Runs(10).Font.Color.RGB = 255
$tr6.Runs(11).Text = ""

# ---------------------------------------------------------------------
# Slide 8 ("JVM vs Dalvik") - Content Placeholder 2 (shape id 12)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

# "Android uses DEX Compiler " (text unchanged, just recolor)
$tr8.Runs(6).Font.Color.RGB = 255

# ---------------------------------------------------------------------
# Slide 9 ("Android SDK") - Content Placeholder 2 (shape id 3)
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9a = $s9.Shapes.Item(2).TextFrame.TextRange

# "Android Software Development Kit " (text unchanged, just recolor)
$tr9a.Runs(1).Font.Color.RGB = 255

# "Development tools" / ": Android Studio IDE and other tools ("
$tr9a.Runs(5).Text = "Development tools: "
$tr9a.Runs(6).Text = "Android Studio IDE and other tools ("

# ---------------------------------------------------------------------
# Slide 9 ("Android SDK") - Content Placeholder 3 (shape id 4)
# ---------------------------------------------------------------------
$tr9b = $s9.Shapes.Item(3).TextFrame.TextRange

# "Android support" / ": extra "
$tr9b.Runs(3).Text = "Android support: "
$tr9b.Runs(3).Font.Color.RGB = 255
$tr9b.Runs(4).Text = "extra "

# "Sample apps:  " (text unchanged, just recolor)
$tr9b.Runs(5).Font.Color.RGB = 255
